# ValueSet-med-busulfan-vs.xlsx update
# "updated with all changes for cibmtr-reporting-ig"
#
# This script edits the "Metadata" worksheet (sheet1) of the CIBMTR
# reporting IG ValueSet workbook:
#   - Version bumped 0.1.6 -> 0.1.7
#   - Status changed active -> draft
#   - Date updated to the new publication timestamp
#   - Contact (organization) now includes the CIBMTR URL
#   - A second Contact row is added for the individual contact (Bob Milius)
#   - A new "Jurisdiction" property row is inserted
#   - Description / Purpose / Copyright / Immutable rows shift down one row
#     to make room for the new Jurisdiction row
#
# The "Include from RxNorm" worksheet (sheet2) content is unchanged.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------
# Make room for the new "Jurisdiction" row by growing the table by one
# row first. We prime row 16 with the same cell formatting (style) as
# the rest of the data rows (copied from row 11, which uses that format)
# so the newly appended row matches the look of the existing table
# instead of picking up a blank default style.
# ---------------------------------------------------------------------
$ws1.Range("A11:B11").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Shift the tail of the property table down one row (bottom-up so we
# never overwrite a value before it has been copied along):
#   old row 15 (Immutable / BooleanType[null])  -> row 16
#   old row 14 (Copyright / "")                 -> row 15
#   old row 13 (Purpose / "")                    -> row 14
#   old row 12 (Description / RxNorm codes...)   -> row 13
# ---------------------------------------------------------------------
$ws1.Range("A16").Value = "Immutable"
$ws1.Range("B16").Value = "BooleanType[null]"

$ws1.Range("A15").Value = "Copyright"
$ws1.Range("B15").Value = ""

$ws1.Range("A14").Value = "Purpose"
$ws1.Range("B14").Value = ""

$ws1.Range("A13").Value = "Description"
$ws1.Range("B13").Value = "RxNorm codes for Busulfan"

# ---------------------------------------------------------------------
# New row 12: Jurisdiction property (no display value).
# ---------------------------------------------------------------------
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

# ---------------------------------------------------------------------
# Row 11: second Contact entry - the individual contact.
# ---------------------------------------------------------------------
$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# ---------------------------------------------------------------------
# Row 10: first Contact entry - now includes the organization URL.
# ---------------------------------------------------------------------
$ws1.Range("A10").Value = "Contact"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# ---------------------------------------------------------------------
# Row 8: publication Date.
# ---------------------------------------------------------------------
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# ---------------------------------------------------------------------
# Row 6: Status.
# ---------------------------------------------------------------------
$ws1.Range("B6").Value = "draft"

# ---------------------------------------------------------------------
# Row 3: Version.
# ---------------------------------------------------------------------
$ws1.Range("B3").Value = "0.1.7"
